# Refresh cryptos list (prices / 1h volume deltas) per the latest crawl.
# Two coins (Hedera / PancakeSwap, and PaxDollar / BabyDogeCoin) also swapped
# rank order, so their whole rows (name, link, price, volume) are rewritten.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.550.76'
$ws.Range("E2").Value = '  +0.37%  '

$ws.Range("D3").Value = '1.817.16'
$ws.Range("E3").Value = '  +1.49%  '

$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").Value = '''228.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.12%  '

$ws.Range("D6").Value = '''0.580'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.23%  '

$ws.Range("E7").Value = '  +0.20%  '

$ws.Range("D8").Value = '''35.03'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.18%  '

$ws.Range("E9").Value = '  +1.61%  '

$ws.Range("E10").Value = '  +0.76%  '

$ws.Range("D11").Value = '''0.0953'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.35%  '

$ws.Range("D12").Value = '2.079.32'
$ws.Range("E12").Value = '  +1.48%  '

$ws.Range("D13").Value = '''11.35'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.58%  '

$ws.Range("D14").Value = '1.825.38'
$ws.Range("E14").Value = '  +1.94%  '

$ws.Range("E15").Value = '  +2.07%  '

$ws.Range("D16").Value = '34.554.21'
$ws.Range("E16").Value = '  +0.48%  '

$ws.Range("D17").Value = '''4.37'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.77%  '

$ws.Range("D18").Value = '''69.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.02%  '

$ws.Range("E19").Value = '  +0.52%  '

$ws.Range("D20").Value = '''245.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.32%  '

$ws.Range("D21").Value = '''11.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.56%  '

$ws.Range("E22").Value = '  +0.22%  '

$ws.Range("E23").Value = '  +0.76%  '

$ws.Range("D24").Value = '''171.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.60%  '

$ws.Range("E25").Value = '  +1.17%  '

$ws.Range("E26").Value = '  +4.51%  '

$ws.Range("D27").Value = '''16.83'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.09%  '

$ws.Range("D28").Value = '''0.119'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.20%  '

$ws.Range("E29").Value = '  +0.00%  '

$ws.Range("D30").Value = '''4.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.73%  '

$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = '''0.0530'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.85%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''1.25'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.28%  '

$ws.Range("E33").Value = '  +1.21%  '

$ws.Range("E34").Value = '  +1.74%  '

$ws.Range("D35").Value = '1.403.28'
$ws.Range("E35").Value = '  -1.02%  '

$ws.Range("D36").Value = '''2.57'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.16%  '

$ws.Range("D37").Value = '''0.681'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.59%  '

$ws.Range("E38").Value = '  +0.96%  '

$ws.Range("E39").Value = '  -0.35%  '

$ws.Range("D40").Value = '''83.35'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.50%  '

$ws.Range("E41").Value = '  +4.37%  '

$ws.Range("E42").Value = '  +2.20%  '

$ws.Range("E43").Value = '  +0.41%  '

$ws.Range("D44").Value = '''13.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.11%  '

$ws.Range("D46").Value = '''0.0510'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.62%  '

$ws.Range("E47").Value = '  +0.28%  '

$ws.Range("D48").Value = '1.979.12'
$ws.Range("E48").Value = '  +1.51%  '

$ws.Range("D49").Value = '''105.62'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.25%  '

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0131'
$ws.Range("E50").Value = '  +1.00%  '

$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").Value = '''1.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.18%  '
